$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 475
$ws.Range("I2").Value = 450
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 450
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -337
$ws.Range("N2").Value = -726
$ws.Range("H18").Value = 2856.5715
$ws.Range("I18").Value = 999
$ws.Range("K18").Value = 999
$ws.Range("M18").Value = -715
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("H29").Value = 3107.4546
$ws.Range("J29").Value = 6561.4
$ws.Range("L29").Value = 19684.2
$ws.Range("N29").Value = -20246.2
$ws.Range("H38").Value = 636.875
$ws.Range("I38").Value = 156.42857
$ws.Range("J38").Value = 4000
$ws.Range("K38").Value = 469.28571
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = -97.28570999999999
$ws.Range("N38").Value = -12744
$ws.Range("H43").Value = 6933.9165
$ws.Range("I43").Value = 6072.45
$ws.Range("J43").Value = 8010.75
$ws.Range("K43").Value = 6072.45
$ws.Range("L43").Value = 8010.75
$ws.Range("M43").Value = -6003.45
$ws.Range("N43").Value = -8148.75
$ws.Range("H51").Value = 8221.5
$ws.Range("J51").Value = 8526
$ws.Range("L51").Value = 8526
$ws.Range("N51").Value = -9494
$ws.Range("H58").Value = 1087.8
$ws.Range("J58").Value = 2174.75
$ws.Range("L58").Value = 6524.25
$ws.Range("N58").Value = -6824.25
$ws.Range("H70").Value = 8473.75
$ws.Range("J70").Value = 9631.666999999999
$ws.Range("L70").Value = 28895.001
$ws.Range("N70").Value = -29435.001
$ws.Range("H73").Value = 8473.75
$ws.Range("J73").Value = 9631.666999999999
$ws.Range("L73").Value = 28895.001
$ws.Range("N73").Value = -30767.001
$ws.Range("H80").Value = 1693.4
$ws.Range("I80").Value = 191.53847
$ws.Range("J80").Value = 3320.4167
$ws.Range("K80").Value = 574.61541
$ws.Range("L80").Value = 9961.250100000001
$ws.Range("M80").Value = 423.38459
$ws.Range("N80").Value = -11957.2501
$ws.Range("H83").Value = 1693.4
$ws.Range("I83").Value = 191.53847
$ws.Range("J83").Value = 3320.4167
$ws.Range("K83").Value = 1723.84623
$ws.Range("L83").Value = 29883.7503
$ws.Range("M83").Value = 3268.15377
$ws.Range("N83").Value = -39867.7503
$ws.Range("H88").Value = 6361.5557
$ws.Range("I88").Value = 5500
$ws.Range("J88").Value = 6607.7144
$ws.Range("K88").Value = 5500
$ws.Range("L88").Value = 6607.7144
$ws.Range("M88").Value = -5094
$ws.Range("N88").Value = -7419.7144
$ws.Range("H91").Value = 6361.5557
$ws.Range("I91").Value = 5500
$ws.Range("J91").Value = 6607.7144
$ws.Range("K91").Value = 5500
$ws.Range("L91").Value = 6607.7144
$ws.Range("M91").Value = -4096
$ws.Range("N91").Value = -9415.714400000001
$ws.Range("H132").Value = 1095.4073
$ws.Range("I132").Value = 1059
$ws.Range("K132").Value = 3177
$ws.Range("M132").Value = -647
$ws.Range("H135").Value = 738.5
$ws.Range("I135").Value = 591.8333
$ws.Range("K135").Value = 5326.4997
$ws.Range("M135").Value = -2791.4997
$ws.Range("H141").Value = 662.4
$ws.Range("I141").Value = 706.1429000000001
$ws.Range("J141").Value = 50
$ws.Range("K141").Value = 2118.4287
$ws.Range("L141").Value = 150
$ws.Range("M141").Value = 3061.5713
$ws.Range("N141").Value = -10510
$ws.Range("N21").ClearContents()
$ws.Range("N23").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18183676
$ws.Range("I32").Value = 18869300
$ws.Range("K32").Value = 18869300
$ws.Range("M32").Value = -18869013
$ws.Range("H132").Value = 2830.8333
$ws.Range("I132").Value = 1914.6136
$ws.Range("K132").Value = 5743.8408
$ws.Range("M132").Value = -3213.8408

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2791.182
$ws.Range("I20").Value = 2055.2104
$ws.Range("J20").Value = 3790
$ws.Range("K20").Value = 2055.2104
$ws.Range("L20").Value = 3790
$ws.Range("M20").Value = -1808.2104
$ws.Range("N20").Value = -4284
$ws.Range("H86").Value = 2700.2307
$ws.Range("I86").Value = 2592.4783
$ws.Range("K86").Value = 2592.4783
$ws.Range("M86").Value = -1469.4783
$ws.Range("H89").Value = 2700.2307
$ws.Range("I89").Value = 2592.4783
$ws.Range("K89").Value = 12962.3915
$ws.Range("M89").Value = -7346.391500000002

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 154432.58
$ws.Range("I31").Value = 9000
$ws.Range("K31").Value = 9000
$ws.Range("M31").Value = -8705
$ws.Range("H34").Value = 154432.58
$ws.Range("I34").Value = 9000
$ws.Range("K34").Value = 9000
$ws.Range("M34").Value = -8798
$ws.Range("H36").Value = 50
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("H40").Value = 50
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("H62").Value = 6315.5
$ws.Range("I62").Value = 3003
$ws.Range("J62").Value = 16253
$ws.Range("K62").Value = 3003
$ws.Range("L62").Value = 16253
$ws.Range("M62").Value = -2379
$ws.Range("N62").Value = -17501
$ws.Range("H65").Value = 6315.5
$ws.Range("I65").Value = 3003
$ws.Range("J65").Value = 16253
$ws.Range("K65").Value = 15015
$ws.Range("L65").Value = 81265
$ws.Range("M65").Value = -11895
$ws.Range("N65").Value = -87505
$ws.Range("N36").ClearContents()
$ws.Range("N40").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1766.6666
$ws.Range("I25").Value = 150
$ws.Range("K25").Value = 450
$ws.Range("M25").Value = -281
$ws.Range("H30").Value = 1766.6666
$ws.Range("I30").Value = 150
$ws.Range("K30").Value = 450
$ws.Range("M30").Value = -348

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 7509.5
$ws.Range("J40").Value = 7509.5
$ws.Range("L40").Value = 7509.5
$ws.Range("N40").Value = -7811.5
$ws.Range("H70").Value = 11374.4
$ws.Range("I70").Value = 7546.5
$ws.Range("J70").Value = 15202.3
$ws.Range("K70").Value = 7546.5
$ws.Range("L70").Value = 15202.3
$ws.Range("M70").Value = -7276.5
$ws.Range("N70").Value = -15742.3
$ws.Range("H73").Value = 11374.4
$ws.Range("I73").Value = 7546.5
$ws.Range("J73").Value = 15202.3
$ws.Range("K73").Value = 7546.5
$ws.Range("L73").Value = 15202.3
$ws.Range("M73").Value = -6610.5
$ws.Range("N73").Value = -17074.3
$ws.Range("H132").Value = 1878338.4
$ws.Range("I132").Value = 1878338.4
$ws.Range("K132").Value = 5635015.199999999
$ws.Range("M132").Value = -5632485.199999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7023.6113
$ws.Range("I22").Value = 3749.5
$ws.Range("J22").Value = 7959.0713
$ws.Range("K22").Value = 3749.5
$ws.Range("L22").Value = 7959.0713
$ws.Range("M22").Value = -3454.5
$ws.Range("N22").Value = -8549.0713
$ws.Range("H27").Value = 7023.6113
$ws.Range("I27").Value = 3749.5
$ws.Range("J27").Value = 7959.0713
$ws.Range("K27").Value = 3749.5
$ws.Range("L27").Value = 7959.0713
$ws.Range("M27").Value = -3642.5
$ws.Range("N27").Value = -8173.0713
$ws.Range("H68").Value = 3808.75
$ws.Range("I68").Value = 2921.4075
$ws.Range("J68").Value = 8600.4
$ws.Range("K68").Value = 2921.4075
$ws.Range("L68").Value = 8600.4
$ws.Range("M68").Value = -2172.4075
$ws.Range("N68").Value = -10098.4
$ws.Range("H71").Value = 3808.75
$ws.Range("I71").Value = 2921.4075
$ws.Range("J71").Value = 8600.4
$ws.Range("K71").Value = 14607.0375
$ws.Range("L71").Value = 43002
$ws.Range("M71").Value = -10863.0375
$ws.Range("N71").Value = -50490
$ws.Range("H132").Value = 3848.772
$ws.Range("J132").Value = 8122.1177
$ws.Range("L132").Value = 24366.3531
$ws.Range("N132").Value = -29426.3531

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 12684.857
$ws.Range("J41").Value = 14558.8
$ws.Range("L41").Value = 14558.8
$ws.Range("N41").Value = -15338.8
$ws.Range("H132").Value = 6646.2964
$ws.Range("I132").Value = 4658.1396
$ws.Range("K132").Value = 13974.4188
$ws.Range("M132").Value = -11444.4188

Write-Output "Applied all Leve profit updates."